$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New query row: "falsdk" with 2 results (semantic search additions)
$ws.Range("A3").Value = "falsdk"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 68.29260699999577
$ws.Range("G3").Value = 0.16801000000123167
$ws.Range("H3").Value = 69.31573699999717
$ws.Range("J3").Value = "Falsk som vatten"
$ws.Range("K3").Value = "One False Move"

# Touch the remaining row-3 cells so they materialize as blank cells
# (matches the source row's empty placeholder cells) without introducing
# any style delta.
$ws.Range("B3").WrapText = $false
$ws.Range("L3:X3").WrapText = $false
